# Updated cryptos list on Tue Jul 11 07:23:57 UTC 2023 with GitHub Actions
# Refresh the latest price/volume figures scraped for the cryptos worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.577.02"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").Value = "1.884.21"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'247.09"
$ws.Range("E5").Value = "  +5.87%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.4771"
$ws.Range("E7").Value = "  +2.06%  "
$ws.Range("D8").Value = "'0.2923"
$ws.Range("E8").Value = "  +3.24%  "
$ws.Range("D9").Value = "'0.06536"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").Value = "'22.05"
$ws.Range("E10").Value = "  +5.21%  "
$ws.Range("D11").Value = "'0.07727"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "'97.73"
$ws.Range("E12").Value = "  +4.50%  "
$ws.Range("D13").Value = "'0.7424"
$ws.Range("E13").Value = "  +9.46%  "
$ws.Range("D14").Value = "1.883.58"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "'5.164"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").Value = "'275.26"
$ws.Range("E16").Value = "  +3.49%  "
$ws.Range("D17").Value = "30.580.94"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "'13.54"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").Value = "'0.000007588"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "2.123.51"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'5.278"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").Value = "'6.217"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("D25").Value = "'9.346"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "'163.71"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("D27").Value = "'18.96"
$ws.Range("E27").Value = "  +2.58%  "
$ws.Range("D28").Value = "'1.951"
$ws.Range("E28").Value = "  +3.66%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.374"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.09997"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").Value = "'1.516"
$ws.Range("E31").Value = "  +4.44%  "
$ws.Range("D32").Value = "'4.333"
$ws.Range("E32").Value = "  +3.20%  "
$ws.Range("D33").Value = "'4.112"
$ws.Range("E33").Value = "  +3.49%  "
$ws.Range("D34").Value = "'0.04812"
$ws.Range("E34").Value = "  +3.52%  "
$ws.Range("D35").Value = "'1.131"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("D36").Value = "'0.7047"
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("D37").Value = "'2.718"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'0.01873"
$ws.Range("E38").Value = "  +3.41%  "
$ws.Range("D39").Value = "'2.751"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").Value = "'6.333"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").Value = "'1.974"
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("D42").Value = "'71.43"
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").Value = "'0.4225"
$ws.Range("E43").Value = "  +4.66%  "
$ws.Range("D44").Value = "'0.8408"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").Value = "'0.9998"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'102.91"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").Value = "'9.303"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").Value = "'7.115"
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("D49").Value = "'35.69"
$ws.Range("E49").Value = "  +4.85%  "
$ws.Range("D50").Value = "'917.75"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "'0.3902"
$ws.Range("E51").Value = "  +4.61%  "
